$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Link Budget")

# Update TX Power PA value (E3): 10 -> -10
$ws.Range("E3").Value = -10

# Update Distance value (E9): 10 -> 2
$ws.Range("E9").Value = 2

# Remove the intermediate "Eb over N0 (lin)" row (old row 27); this shifts
# "Bit error rate" row (old row 28) and the trailing blank row (old row 29) up.
$ws.Rows("27").Delete()

# The Eb/N0 row (now row 26) re-wraps to two lines once the sheet settles.
$ws.Rows("26").RowHeight = 34

$ws.Range("J12").Select()
